# Fix port data: several "region totals" rows were missing the port name
# in column B (it had slid into column C, the species column), and several
# "Total check" rows were missing their 0 values. Also distinguish rows
# where "All other" was actually the sole species ("All species") or one
# of several partial species ("All other species"). Finally, add an
# AutoFilter over the data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Region-totals rows: port name had been typed into column C (species)
# instead of column B (port); move it over and set C to "Totals". ---
function Fix-RegionTotalsRow($row, $regionText) {
    $ws.Range("B$row").Value = $regionText
    $ws.Range("C$row").Value = "Totals"
}

Fix-RegionTotalsRow 2   "Eureka region totals"
Fix-RegionTotalsRow 58  "Sacramento region totals"
Fix-RegionTotalsRow 86  "San Francisco region totals"
Fix-RegionTotalsRow 157 "Monterey region totals"
Fix-RegionTotalsRow 192 "Santa Barbara region totals"
Fix-RegionTotalsRow 257 "Los Angeles region totals"
Fix-RegionTotalsRow 367 "San Diego region totals"

# --- "All other" species rows that are really one of several partial
# categories ("All other species"), or the single species for that port
# ("All species"). (Order matters for shared-string creation order.) ---
$ws.Range("C364").Value = "All other species"
$ws.Range("C55").Value = "All species"
$ws.Range("C154").Value = "All species"
$ws.Range("C400").Value = "All species"

# --- "Total check" rows that were missing their computed 0 values. ---
$checkRows = @(57, 85, 156, 250, 253, 256, 395, 402)
foreach ($r in $checkRows) {
    $ws.Range("D$r").Value = 0
    $ws.Range("E$r").Value = 0
}

# --- Add an AutoFilter over the full data range (this also registers the
# sheet-scoped hidden "_FilterDatabase" defined name that Excel creates
# behind the scenes). ---
$ws.Range("A1:E402").AutoFilter() | Out-Null
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$E`$402")
$filterName.Visible = $false

# --- Restore the user's last cell selection. ---
$ws.Range("C10").Select() | Out-Null
